{"js": "// \"fixed date in CV\" \u2014 the PwC summer-placement dates were a year out of\n// date (2021) and need to read 2022 instead, e.g.\n//   13/06/2021 \u2013 19/08/2021   ->   13/06/2022 \u2013 19/08/2022\n//\n// Use Word's search API to locate each date string and replace it in\n// place so surrounding run formatting (font size, style, etc.) is left\n// untouched.\n\nconst replacements = [\n  { find: \"13/06/2021\", replace: \"13/06/2022\" },\n  { find: \"19/08/2021\", replace: \"19/08/2022\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# \"fixed date in CV\" \u2014 the PwC summer-placement dates were a year out of\n# date (2021) and need to read 2022 instead, e.g.\n#   13/06/2021 - 19/08/2021   ->   13/06/2022 - 19/08/2022\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @{ Find = \"13/06/2021\"; Replace = \"13/06/2022\" },\n    @{ Find = \"19/08/2021\"; Replace = \"19/08/2022\" }\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $range.Find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.Replace, $wdReplaceAll) | Out-Null\n}\n"}
